$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44727
$ws.Range("J2").Value2 = 400
$ws.Range("K2").Value2 = 18000
$ws.Range("L2").Value2 = 19000
$ws.Range("M2").Value2 = 18500
$ws.Range("P2").Value2 = 1850

$ws.Range("D3").Value2 = 44644
$ws.Range("J3").Value2 = 300
$ws.Range("K3").Value2 = 20000
$ws.Range("L3").Value2 = 21000
$ws.Range("M3").Value2 = 20500
$ws.Range("P3").Value2 = 2050

$ws.Range("D4").Value2 = 44524
$ws.Range("J4").Value2 = 200
$ws.Range("K4").Value2 = 20000
$ws.Range("L4").Value2 = 21000
$ws.Range("M4").Value2 = 20500
$ws.Range("P4").Value2 = 2050

$ws.Range("D5").Value2 = 44291
$ws.Range("J5").Value2 = 200
$ws.Range("K5").Value2 = 13000
$ws.Range("L5").Value2 = 14000
$ws.Range("M5").Value2 = 13500
$ws.Range("P5").Value2 = 1350

$ws.Range("D6").Value2 = 44428
$ws.Range("J6").Value2 = 300
$ws.Range("K6").Value2 = 15000
$ws.Range("L6").Value2 = 16000
$ws.Range("M6").Value2 = 15500
$ws.Range("P6").Value2 = 1550

$ws.Range("D8").Value2 = 44377
$ws.Range("J8").Value2 = 650
$ws.Range("K8").Value2 = 14000
$ws.Range("L8").Value2 = 15000
$ws.Range("M8").Value2 = 14538
$ws.Range("P8").Value2 = 1454

$ws.Range("D9").Value2 = 44777
$ws.Range("J9").Value2 = 200
$ws.Range("K9").Value2 = 24000
$ws.Range("L9").Value2 = 25000
$ws.Range("M9").Value2 = 24500
$ws.Range("P9").Value2 = 2450

$ws.Range("D10").Value2 = 44714
$ws.Range("J10").Value2 = 400
$ws.Range("K10").Value2 = 19000
$ws.Range("L10").Value2 = 20000
$ws.Range("M10").Value2 = 19500
$ws.Range("P10").Value2 = 1950

$ws.Range("D11").Value2 = 44441
$ws.Range("J11").Value2 = 300
$ws.Range("K11").Value2 = 15000
$ws.Range("L11").Value2 = 16000
$ws.Range("M11").Value2 = 15500
$ws.Range("P11").Value2 = 1550

$ws.Range("D12").Value2 = 44358
$ws.Range("J12").Value2 = 300
$ws.Range("K12").Value2 = 14000
$ws.Range("L12").Value2 = 15000
$ws.Range("M12").Value2 = 14500
$ws.Range("P12").Value2 = 1450

$ws.Range("D13").Value2 = 44460
$ws.Range("J13").Value2 = 300
$ws.Range("K13").Value2 = 15000
$ws.Range("L13").Value2 = 16000
$ws.Range("M13").Value2 = 15500
$ws.Range("P13").Value2 = 1550

$ws.Range("D14").Value2 = 44160
$ws.Range("J14").Value2 = 360
$ws.Range("K14").Value2 = 10000
$ws.Range("L14").Value2 = 11000
$ws.Range("M14").Value2 = 10500
$ws.Range("P14").Value2 = 1050

$ws.Range("D15").Value2 = 44218
$ws.Range("J15").Value2 = 320
$ws.Range("K15").Value2 = 10000
$ws.Range("L15").Value2 = 11000
$ws.Range("M15").Value2 = 10500
$ws.Range("P15").Value2 = 1050

$ws.Range("D16").Value2 = 44263
$ws.Range("J16").Value2 = 300
$ws.Range("K16").Value2 = 15000
$ws.Range("L16").Value2 = 16000
$ws.Range("M16").Value2 = 15500
$ws.Range("P16").Value2 = 1550

$ws.Range("D17").Value2 = 44265
$ws.Range("J17").Value2 = 200
$ws.Range("K17").Value2 = 15000
$ws.Range("L17").Value2 = 16000
$ws.Range("M17").Value2 = 15500
$ws.Range("P17").Value2 = 1550

$ws.Range("D18").Value2 = 44694
$ws.Range("J18").Value2 = 400
$ws.Range("K18").Value2 = 16000
$ws.Range("L18").Value2 = 17000
$ws.Range("M18").Value2 = 16500
$ws.Range("P18").Value2 = 1650

$ws.Range("D19").Value2 = 44406
$ws.Range("J19").Value2 = 400
$ws.Range("K19").Value2 = 14000
$ws.Range("L19").Value2 = 15000
$ws.Range("M19").Value2 = 14500
$ws.Range("P19").Value2 = 1450

$ws.Range("D20").Value2 = 44547
$ws.Range("J20").Value2 = 300
$ws.Range("K20").Value2 = 19000
$ws.Range("L20").Value2 = 20000
$ws.Range("M20").Value2 = 19500
$ws.Range("P20").Value2 = 1950

$ws.Range("D21").Value2 = 44679
$ws.Range("J21").Value2 = 200
$ws.Range("K21").Value2 = 19000
$ws.Range("L21").Value2 = 20000
$ws.Range("M21").Value2 = 19500
$ws.Range("P21").Value2 = 1950

$ws.Range("D22").Value2 = 44330
$ws.Range("J22").Value2 = 300
$ws.Range("K22").Value2 = 13000
$ws.Range("L22").Value2 = 14000
$ws.Range("M22").Value2 = 13500
$ws.Range("P22").Value2 = 1350

$ws.Range("D23").Value2 = 44204
$ws.Range("J23").Value2 = 400
$ws.Range("K23").Value2 = 10000
$ws.Range("L23").Value2 = 11000
$ws.Range("M23").Value2 = 10500
$ws.Range("P23").Value2 = 1050
